$d = $word.ActiveDocument

# --- Change 1: add a new bold heading line above the existing first paragraph ---
# Emulates placing the cursor at the very start of paragraph 1, typing the new
# heading text, then pressing Enter to push the old text into its own paragraph.
$firstPara = $d.Paragraphs.Item(1)
$startPos = $firstPara.Range.Start
$headingText = "Git Guide " + [char]0x2014 + " update check from Nadira "

$insertPoint = $d.Range($startPos, $startPos)
$insertPoint.InsertBefore($headingText)

$breakPos = $startPos + $headingText.Length
$breakPoint = $d.Range($breakPos, $breakPos)
$breakPoint.InsertParagraphAfter()

# --- Change 2: append a blank paragraph and a NOTE paragraph after "Git push" ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$blankPara.Range.InsertParagraphAfter()

$noteText = "NOTE: We might want to fix the " + [char]0x201C + "detached HEAD" + [char]0x201D + " "
$notePara = $d.Paragraphs.Item($d.Paragraphs.Count)
$notePara.Range.Text = $noteText
